$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.773.75'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '2.424.73'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.74%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.35%  '
$ws.Range('E9').Value = '  +4.89%  '
$ws.Range('E10').Value = '  +3.08%  '
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.70'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('D14').Value = '2.853.30'
$ws.Range('E14').Value = '  +2.83%  '
$ws.Range('D15').Value = '59.659.98'
$ws.Range('E15').Value = '  +2.28%  '
$ws.Range('E16').Value = '  +4.44%  '
$ws.Range('D17').Value = '2.405.87'
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('E18').Value = '  +5.72%  '
$ws.Range('E19').Value = '  +4.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '334.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.24%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.85%  '
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('D28').Value = '0.0₃0787'
$ws.Range('E28').Value = '  +6.72%  '
$ws.Range('E29').Value = '  +2.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.61'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  +2.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.73%  '
$ws.Range('E33').Value = '  +0.72%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  +5.23%  '
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '40.08'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.29%  '
$ws.Range('E40').Value = '  +11.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '315.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.82%  '
$ws.Range('E42').Value = '  +2.79%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0968'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.77%  '
$ws.Range('E45').Value = '  +4.26%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.574'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.403'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.54%  '
$ws.Range('E49').Value = '  +3.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.05'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('E51').Value = '  +5.01%  '
